# Insert a new weekly price-report row before row 12 (shifts old rows 12-43 -> 13-44)
# and populate the newly inserted row 12 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; existing rows 12..43 shift down to 13..44.
$ws.Rows.Item(12).Insert()

# Fill in the new row 12 with the data for this entry.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44487
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112012
$ws.Cells.Item(12, 7).Value = "Espinaca"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 950
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = 975
$ws.Cells.Item(12, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 325
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = "Hortaliza"
